$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before column E ("Padre rubro"), shifting old E:Z to F:AA
$ws.Columns.Item(5).EntireColumn.Insert()
$ws.Cells.Item(1,5).Value2 = "Padre rubro"
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# 2) Append two new trailing columns after the (shifted) last column Z: AA already holds old "m3" (shifted),
#    so we add a brand-new column AB = "Codigo Prov"
$ws.Cells.Item(1,28).Value2 = "Código Prov"

# Give AB2 the same style as the other "list" formatted cells (e.g. Q2 which now uses style 3)
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("AB2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 3) Update the view to reflect scrolling to the new last column and selecting AB2
$ws.Range("AB2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 22   # column V
$excel.ActiveWindow.ScrollRow = 1

Write-Host "Done"
